$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.160.77"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "'2.297.95"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'313.11"
$ws.Range("E5").Value = "  -3.81%  "
$ws.Range("D6").Value = "'106.57"
$ws.Range("E6").Value = "  +3.25%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -2.27%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").Value = "'40.28"
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "'8.30"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "'0.973"
$ws.Range("E14").Value = "  -3.55%  "
$ws.Range("D15").Value = "'15.56"
$ws.Range("E15").Value = "  -3.74%  "
$ws.Range("D16").Value = "'2.644.69"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "'2.299.61"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "'41.977.72"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").Value = "'7.54"
$ws.Range("E19").Value = "  -4.74%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "'73.25"
$ws.Range("E21").Value = "  -4.91%  "
$ws.Range("E22").Value = "  -5.51%  "
$ws.Range("D23").Value = "'257.30"
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -0.63%  "
$ws.Range("D25").Value = "'9.38"
$ws.Range("E25").Value = "  -5.64%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.28"
$ws.Range("E28").Value = "  +2.90%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'22.85"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "'166.49"
$ws.Range("E30").Value = "  -4.72%  "
$ws.Range("D31").Value = "'35.70"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("E33").Value = "  -7.02%  "
$ws.Range("D34").Value = "'5.80"
$ws.Range("E34").Value = "  -7.66%  "
$ws.Range("E35").Value = "  +5.61%  "
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").Value = "'4.61"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").Value = "'2.89"
$ws.Range("E39").Value = "  +6.09%  "
$ws.Range("D40").Value = "'3.63"
$ws.Range("E40").Value = "  -4.59%  "
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "'71.88"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "'97.26"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("D47").Value = "'113.36"
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("D48").Value = "'9.16"
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  -4.87%  "
$ws.Range("D50").Value = "'75.75"
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("E51").Value = "  -0.85%  "
